# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.103.78"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "2.378.08"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").Value = "2.376.33"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.338"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").Value = "2.794.32"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("E16").Value = "  -1.87%  "

$ws.Range("D17").Value = "59.905.94"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").Value = "2.364.85"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.52%  "

$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "556.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.79%  "

$ws.Range("D29").Value = "2.494.67"
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("D30").Value = "0.0₃0928"
$ws.Range("E30").Value = "  +2.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.85%  "

$ws.Range("E32").Value = "  -2.33%  "

$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("E36").Value = "  +5.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.71%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("E39").Value = "  -0.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("E45").Value = "  +4.14%  "

$ws.Range("D46").Value = "0.0₆0291"
$ws.Range("E46").Value = "  +3.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.587"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0501"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
